$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.268639
$ws.Range("H2").Value = 60.80591700000001
$ws.Range("I2").Value = 0.09537690375401175
$ws.Range("J2").Value = 0.09537690375401174
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.002302
$ws.Range("N2").Value = 0.006906
$ws.Range("O2").Value = 0.002384587814362636
$ws.Range("P2").Value = 0.002384587814362636
$ws.Range("Q2").Value = 0.04665840697800001
$ws.Range("R2").Value = 0.4199256628020001
$ws.Range("S2").Value = 0.0002274346024634543
$ws.Range("T2").Value = 0.0002274346024634543
$ws.Range("G3").Value = 20.268639
$ws.Range("H3").Value = 60.80591700000001
$ws.Range("I3").Value = 0.09537690375401175
$ws.Range("J3").Value = 0.09537690375401174
$ws.Range("M3").Value = 0.01379866666666667
$ws.Range("N3").Value = 0.041396
$ws.Range("O3").Value = 0.01429371519886413
$ws.Range("P3").Value = 0.01429371519886413
$ws.Range("Q3").Value = 0.2796801933480001
$ws.Range("R3").Value = 2.517121740132001
$ws.Range("S3").Value = 0.001363290298809319
$ws.Range("T3").Value = 0.001363290298809319
$ws.Range("G4").Value = 20.268639
$ws.Range("H4").Value = 60.80591700000001
$ws.Range("I4").Value = 0.09537690375401175
$ws.Range("J4").Value = 0.09537690375401174
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.9492653333333333
$ws.Range("N4").Value = 2.847796
$ws.Range("O4").Value = 0.9833216969867733
$ws.Range("P4").Value = 0.9833216969867732
$ws.Range("Q4").Value = 19.240316356548
$ws.Range("R4").Value = 173.162847208932
$ws.Range("S4").Value = 0.09378617885273899
$ws.Range("T4").Value = 0.09378617885273896
$ws.Range("I5").Value = 0.6011031624655011
$ws.Range("J5").Value = 0.601103162465501
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.002302
$ws.Range("N5").Value = 0.006906
$ws.Range("O5").Value = 0.002384587814362636
$ws.Range("P5").Value = 0.002384587814362636
$ws.Range("Q5").Value = 0.2940598288073333
$ws.Range("R5").Value = 2.646538459266
$ws.Range("S5").Value = 0.001433383276390078
$ws.Range("T5").Value = 0.001433383276390077
$ws.Range("I6").Value = 0.6011031624655011
$ws.Range("J6").Value = 0.601103162465501
$ws.Range("M6").Value = 0.01379866666666667
$ws.Range("N6").Value = 0.041396
$ws.Range("O6").Value = 0.01429371519886413
$ws.Range("P6").Value = 0.01429371519886413
$ws.Range("Q6").Value = 1.762655759239556
$ws.Range("R6").Value = 15.863901833156
$ws.Range("S6").Value = 0.008591997409418427
$ws.Range("T6").Value = 0.008591997409418425
$ws.Range("I7").Value = 0.6011031624655011
$ws.Range("J7").Value = 0.601103162465501
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9492653333333333
$ws.Range("N7").Value = 2.847796
$ws.Range("O7").Value = 0.9833216969867733
$ws.Range("P7").Value = 0.9833216969867732
$ws.Range("Q7").Value = 121.2601222470618
$ws.Range("R7").Value = 1091.341100223556
$ws.Range("S7").Value = 0.5910777817796926
$ws.Range("T7").Value = 0.5910777817796925
$ws.Range("G8").Value = 42.02733833333333
$ws.Range("H8").Value = 126.082015
$ws.Range("I8").Value = 0.1977654939365007
$ws.Range("J8").Value = 0.1977654939365007
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.002302
$ws.Range("N8").Value = 0.006906
$ws.Range("O8").Value = 0.002384587814362636
$ws.Range("P8").Value = 0.002384587814362636
$ws.Range("Q8").Value = 0.09674693284333331
$ws.Range("R8").Value = 0.8707223955899999
$ws.Range("S8").Value = 0.0004715891869423872
$ws.Range("T8").Value = 0.0004715891869423872
$ws.Range("G9").Value = 42.02733833333333
$ws.Range("H9").Value = 126.082015
$ws.Range("I9").Value = 0.1977654939365007
$ws.Range("J9").Value = 0.1977654939365007
$ws.Range("M9").Value = 0.01379866666666667
$ws.Range("N9").Value = 0.041396
$ws.Range("O9").Value = 0.01429371519886413
$ws.Range("P9").Value = 0.01429371519886413
$ws.Range("Q9").Value = 0.5799212325488888
$ws.Range("R9").Value = 5.21929109294
$ws.Range("S9").Value = 0.002826803646491031
$ws.Range("T9").Value = 0.002826803646491031
$ws.Range("G10").Value = 42.02733833333333
$ws.Range("H10").Value = 126.082015
$ws.Range("I10").Value = 0.1977654939365007
$ws.Range("J10").Value = 0.1977654939365007
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.9492653333333333
$ws.Range("N10").Value = 2.847796
$ws.Range("O10").Value = 0.9833216969867733
$ws.Range("P10").Value = 0.9833216969867732
$ws.Range("Q10").Value = 39.89509533210443
$ws.Range("R10").Value = 359.0558579889399
$ws.Range("S10").Value = 0.1944671011030673
$ws.Range("T10").Value = 0.1944671011030673
$ws.Range("G11").Value = 22.47397933333333
$ws.Range("H11").Value = 67.421938
$ws.Range("I11").Value = 0.1057544398439867
$ws.Range("J11").Value = 0.1057544398439867
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.002302
$ws.Range("N11").Value = 0.006906
$ws.Range("O11").Value = 0.002384587814362636
$ws.Range("P11").Value = 0.002384587814362636
$ws.Range("Q11").Value = 0.05173510042533332
$ws.Range("R11").Value = 0.4656159038279999
$ws.Range("S11").Value = 0.000252180748566717
$ws.Range("T11").Value = 0.000252180748566717
$ws.Range("G12").Value = 22.47397933333333
$ws.Range("H12").Value = 67.421938
$ws.Range("I12").Value = 0.1057544398439867
$ws.Range("J12").Value = 0.1057544398439867
$ws.Range("M12").Value = 0.01379866666666667
$ws.Range("N12").Value = 0.041396
$ws.Range("O12").Value = 0.01429371519886413
$ws.Range("P12").Value = 0.01429371519886413
$ws.Range("Q12").Value = 0.3101109494942222
$ws.Range("R12").Value = 2.790998545448
$ws.Range("S12").Value = 0.001511623844145355
$ws.Range("T12").Value = 0.001511623844145354
$ws.Range("G13").Value = 22.47397933333333
$ws.Range("H13").Value = 67.421938
$ws.Range("I13").Value = 0.1057544398439867
$ws.Range("J13").Value = 0.1057544398439867
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.9492653333333333
$ws.Range("N13").Value = 2.847796
$ws.Range("O13").Value = 0.9833216969867733
$ws.Range("P13").Value = 0.9833216969867732
$ws.Range("Q13").Value = 21.33376948318311
$ws.Range("R13").Value = 192.003925348648
$ws.Range("S13").Value = 0.1039906352512746
$ws.Range("T13").Value = 0.1039906352512746
